$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Col4a5", "Cd93", "ECs", 2, 0.6666666666666666, 0.126664, 0.379992, 0.02685339687075019, 0.02685339687075019, 3, 1, 229.5846506666667, 688.753952, 0.5033187360873315, 0.5033187360873315, 29.08011019204267, 261.720991728384, 0.01351581777263749, 0.01351581777263748),
    @("ECs", "Col4a5", "Cd93", "M1", 2, 0.6666666666666666, 0.126664, 0.379992, 0.02685339687075019, 0.02685339687075019, 3, 1, 135.7283196666666, 407.1849589999999, 0.2975573763642838, 0.2975573763642838, 17.19189188225866, 154.727026940328, 0.007990426319329294, 0.007990426319329292),
    @("ECs", "Col4a5", "Cd93", "M2", 2, 0.6666666666666666, 0.126664, 0.379992, 0.02685339687075019, 0.02685339687075019, 3, 1, 90.23148833333335, 270.694465, 0.1978146123067711, 0.1978146123067711, 11.42908123825334, 102.86173114428, 0.005311994291107308, 0.005311994291107307),
    @("ECs", "Col4a5", "Cd93", "sCs", 2, 0.6666666666666666, 0.126664, 0.379992, 0.02685339687075019, 0.02685339687075019, 3, 1, 0.5972149999999999, 1.791645, 0.001309275241613694, 0.001309275241613694, 0.07564564075999999, 0.68081076684, 0.00003515848767609988, 0.00003515848767609988),
    @("FAPs", "Col4a5", "Cd93", "ECs", 3, 1, 3.579479333333333, 10.738438, 0.758867390329125, 0.758867390329125, 3, 1, 229.5846506666667, 688.753952, 0.5033187360873315, 0.5033187360873315, 821.7935123118862, 7396.141610806975, 0.3819521757583468, 0.3819521757583468),
    @("FAPs", "Col4a5", "Cd93", "M1", 3, 1, 3.579479333333333, 10.738438, 0.758867390329125, 0.758867390329125, 3, 1, 135.7283196666666, 407.1849589999999, 0.2975573763642838, 0.2975573763642838, 485.8367151948934, 4372.530436754041, 0.2258065896747453, 0.2258065896747453),
    @("FAPs", "Col4a5", "Cd93", "M2", 3, 1, 3.579479333333333, 10.738438, 0.758867390329125, 0.758867390329125, 3, 1, 90.23148833333335, 270.694465, 0.1978146123067711, 0.1978146123067711, 322.9817477050744, 2906.83572934567, 0.150115058610207, 0.150115058610207),
    @("FAPs", "Col4a5", "Cd93", "sCs", 3, 1, 3.579479333333333, 10.738438, 0.758867390329125, 0.758867390329125, 3, 1, 0.5972149999999999, 1.791645, 0.001309275241613694, 0.001309275241613694, 2.137718750056666, 19.23946875051, 0.0009935662858259189, 0.0009935662858259189),
    @("sCs", "Col4a5", "Cd93", "ECs", 3, 1, 1.010727333333333, 3.032182, 0.2142792128001249, 0.2142792128001248, 3, 1, 229.5846506666667, 688.753952, 0.5033187360873315, 0.5033187360873315, 232.0474817425849, 2088.427335683264, 0.1078507425563472, 0.1078507425563472),
    @("sCs", "Col4a5", "Cd93", "M1", 3, 1, 1.010727333333333, 3.032182, 0.2142792128001249, 0.2142792128001248, 3, 1, 135.7283196666666, 407.1849589999999, 0.2975573763642838, 0.2975573763642838, 137.1843225945042, 1234.658903350538, 0.0637603603702092, 0.06376036037020918),
    @("sCs", "Col4a5", "Cd93", "M2", 3, 1, 1.010727333333333, 3.032182, 0.2142792128001249, 0.2142792128001248, 3, 1, 90.23148833333335, 270.694465, 0.1978146123067711, 0.1978146123067711, 91.19943158584779, 820.79488427263, 0.04238755940545679, 0.04238755940545678),
    @("sCs", "Col4a5", "Cd93", "sCs", 3, 1, 1.010727333333333, 3.032182, 0.2142792128001249, 0.2142792128001248, 3, 1, 0.5972149999999999, 1.791645, 0.001309275241613694, 0.001309275241613694, 0.6036215243766666, 5.432593719389999, 0.0002805504681116757, 0.0002805504681116757)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
